$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 625.2857
$ws.Range("I2").Value = 625.2857
$ws.Range("K2").Value = 625.2857
$ws.Range("M2").Value = -512.2857

$ws.Range("H4").Value = 3105.0588
$ws.Range("I4").Value = 218.7
$ws.Range("J4").Value = 7228.4287
$ws.Range("K4").Value = 218.7
$ws.Range("L4").Value = 7228.4287
$ws.Range("M4").Value = -104.7
$ws.Range("N4").Value = -7456.4287

$ws.Range("H7").Value = 30000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 30000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 30000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -30224

$ws.Range("H14").Value = 30000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 30000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 30000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -30382

$ws.Range("H28").Value = 603
$ws.Range("I28").Value = 360.7143
$ws.Range("J28").Value = 885.6667
$ws.Range("K28").Value = 360.7143
$ws.Range("L28").Value = 885.6667
$ws.Range("M28").Value = 124.2857
$ws.Range("N28").Value = -1855.6667

$ws.Range("H32").Value = 952.8570999999999
$ws.Range("I32").Value = 500
$ws.Range("J32").Value = 1076.3636
$ws.Range("K32").Value = 500
$ws.Range("L32").Value = 1076.3636
$ws.Range("M32").Value = -174
$ws.Range("N32").Value = -1728.3636

$ws.Range("H40").Value = 1936.841
$ws.Range("I40").Value = 1954.8572
$ws.Range("J40").Value = 1866.7778
$ws.Range("K40").Value = 1954.8572
$ws.Range("L40").Value = 1866.7778
$ws.Range("M40").Value = -1779.8572
$ws.Range("N40").Value = -2216.7778

$ws.Range("H62").Value = 1243.1111
$ws.Range("I62").Value = 1276.3334
$ws.Range("J62").Value = 1176.6666
$ws.Range("K62").Value = 1276.3334
$ws.Range("L62").Value = 1176.6666
$ws.Range("M62").Value = -652.3334
$ws.Range("N62").Value = -2424.6666

$ws.Range("H65").Value = 1243.1111
$ws.Range("I65").Value = 1276.3334
$ws.Range("J65").Value = 1176.6666
$ws.Range("K65").Value = 6381.666999999999
$ws.Range("L65").Value = 5883.333000000001
$ws.Range("M65").Value = -3261.666999999999
$ws.Range("N65").Value = -12123.333

$ws.Range("H99").Value = 573.1429000000001
$ws.Range("I99").Value = 502
$ws.Range("K99").Value = 1506
$ws.Range("M99").Value = -8

$ws.Range("H107").Value = 8067030.5
$ws.Range("I107").Value = 8621598
$ws.Range("K107").Value = 8621598
$ws.Range("M107").Value = -8619678

$ws.Range("H116").Value = 5588.615
$ws.Range("I116").Value = 8193.6
$ws.Range("J116").Value = 2036.3636
$ws.Range("K116").Value = 8193.6
$ws.Range("L116").Value = 2036.3636
$ws.Range("M116").Value = -4751.6
$ws.Range("N116").Value = -8920.363600000001

$ws.Range("H118").Value = 638.8333
$ws.Range("I118").Value = 261.33334
$ws.Range("J118").Value = 1016.3333
$ws.Range("K118").Value = 784.0000200000001
$ws.Range("L118").Value = 3048.9999
$ws.Range("M118").Value = 872.9999799999999
$ws.Range("N118").Value = -6362.9999

$ws.Range("H132").Value = 1704.7273
$ws.Range("I132").Value = 1704.7273
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5114.1819
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2584.1819
$ws.Range("N132").ClearContents()

$ws.Range("H138").Value = 3453.3704
$ws.Range("I138").Value = 1370.4839
$ws.Range("J138").Value = 4744.76
$ws.Range("K138").Value = 4111.4517
$ws.Range("L138").Value = 14234.28
$ws.Range("M138").Value = 1028.5483
$ws.Range("N138").Value = -24514.28

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7827.1177
$ws.Range("I32").Value = 8703.6875
$ws.Range("J32").Value = 6350.7896
$ws.Range("K32").Value = 8703.6875
$ws.Range("L32").Value = 6350.7896
$ws.Range("M32").Value = -8416.6875
$ws.Range("N32").Value = -6924.7896

$ws.Range("H132").Value = 2078.3174
$ws.Range("I132").Value = 1232.279
$ws.Range("J132").Value = 3897.3
$ws.Range("K132").Value = 3696.837
$ws.Range("L132").Value = 11691.9
$ws.Range("M132").Value = -1166.837
$ws.Range("N132").Value = -16751.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 321.2857
$ws.Range("I22").Value = 321.2857
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 321.2857
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -148.2857
$ws.Range("N22").ClearContents()

$ws.Range("H86").Value = 2702.3333
$ws.Range("J86").Value = 2702.3333
$ws.Range("L86").Value = 2702.3333
$ws.Range("N86").Value = -4948.3333

$ws.Range("H89").Value = 2702.3333
$ws.Range("J89").Value = 2702.3333
$ws.Range("L89").Value = 13511.6665
$ws.Range("N89").Value = -24743.6665

$ws.Range("H107").Value = 1238
$ws.Range("I107").Value = 1216.6428
$ws.Range("J107").Value = 1337.6666
$ws.Range("K107").Value = 1216.6428
$ws.Range("L107").Value = 1337.6666
$ws.Range("M107").Value = 703.3571999999999
$ws.Range("N107").Value = -5177.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1739.46
$ws.Range("I31").Value = 1123.3438
$ws.Range("J31").Value = 2029.3971
$ws.Range("K31").Value = 1123.3438
$ws.Range("L31").Value = 2029.3971
$ws.Range("M31").Value = -828.3438000000001
$ws.Range("N31").Value = -2619.3971

$ws.Range("H34").Value = 1739.46
$ws.Range("I34").Value = 1123.3438
$ws.Range("J34").Value = 2029.3971
$ws.Range("K34").Value = 1123.3438
$ws.Range("L34").Value = 2029.3971
$ws.Range("M34").Value = -921.3438000000001
$ws.Range("N34").Value = -2433.3971

$ws.Range("H62").Value = 3745.2632
$ws.Range("I62").Value = 3597.647
$ws.Range("K62").Value = 3597.647
$ws.Range("M62").Value = -2973.647

$ws.Range("H65").Value = 3745.2632
$ws.Range("I65").Value = 3597.647
$ws.Range("K65").Value = 17988.235
$ws.Range("M65").Value = -14868.235

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 3448.7144
$ws.Range("I24").Value = 129
$ws.Range("J24").Value = 4002
$ws.Range("K24").Value = 387
$ws.Range("L24").Value = 12006
$ws.Range("M24").Value = -157
$ws.Range("N24").Value = -12466

$ws.Range("H68").Value = 2453.719
$ws.Range("I68").Value = 3137.6047
$ws.Range("J68").Value = 1814.4348
$ws.Range("K68").Value = 9412.8141
$ws.Range("L68").Value = 5443.3044
$ws.Range("M68").Value = -8601.8141
$ws.Range("N68").Value = -7065.3044

$ws.Range("H71").Value = 2453.719
$ws.Range("I71").Value = 3137.6047
$ws.Range("J71").Value = 1814.4348
$ws.Range("K71").Value = 28238.4423
$ws.Range("L71").Value = 16329.9132
$ws.Range("M71").Value = -24182.4423
$ws.Range("N71").Value = -24441.9132

$ws.Range("H107").Value = 1260.5238
$ws.Range("I107").Value = 377.5
$ws.Range("J107").Value = 1407.6945
$ws.Range("K107").Value = 1132.5
$ws.Range("L107").Value = 4223.083500000001
$ws.Range("M107").Value = 787.5
$ws.Range("N107").Value = -8063.083500000001

$ws.Range("H121").Value = 1043.0204
$ws.Range("I121").Value = 725
$ws.Range("J121").Value = 1056.5532
$ws.Range("K121").Value = 2175
$ws.Range("L121").Value = 3169.6596
$ws.Range("M121").Value = -865
$ws.Range("N121").Value = -5789.6596

$ws.Range("H122").Value = 517.67346
$ws.Range("I122").Value = 462.66666
$ws.Range("J122").Value = 535.5135
$ws.Range("K122").Value = 4163.99994
$ws.Range("L122").Value = 4819.6215
$ws.Range("M122").Value = -1713.99994
$ws.Range("N122").Value = -9719.621500000001

$ws.Range("H129").Value = 15153598
$ws.Range("I129").Value = 41667956
$ws.Range("J129").Value = 2536.3572
$ws.Range("K129").Value = 125003868
$ws.Range("L129").Value = 7609.071599999999
$ws.Range("M129").Value = -124998868
$ws.Range("N129").Value = -17609.0716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3864
$ws.Range("I80").Value = 3285.5715
$ws.Range("J80").Value = 5888.5
$ws.Range("K80").Value = 3285.5715
$ws.Range("L80").Value = 5888.5
$ws.Range("M80").Value = -2287.5715
$ws.Range("N80").Value = -7884.5

$ws.Range("H83").Value = 3864
$ws.Range("I83").Value = 3285.5715
$ws.Range("J83").Value = 5888.5
$ws.Range("K83").Value = 16427.8575
$ws.Range("L83").Value = 29442.5
$ws.Range("M83").Value = -11435.8575
$ws.Range("N83").Value = -39426.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 19231136
$ws.Range("I55").Value = 301.75
$ws.Range("J55").Value = 35714708
$ws.Range("K55").Value = 301.75
$ws.Range("L55").Value = 35714708
$ws.Range("M55").Value = -128.75
$ws.Range("N55").Value = -35715054

$ws.Range("H122").Value = 11642082
$ws.Range("I122").Value = 17865644
$ws.Range("J122").Value = 3344000
$ws.Range("K122").Value = 53596932
$ws.Range("L122").Value = 10032000
$ws.Range("M122").Value = -53594482
$ws.Range("N122").Value = -10036900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2230.4055
$ws.Range("I132").Value = 1764.2
$ws.Range("J132").Value = 2778.8823
$ws.Range("K132").Value = 5292.6
$ws.Range("L132").Value = 8336.6469
$ws.Range("M132").Value = -2762.6
$ws.Range("N132").Value = -13396.6469
